# Auto-update: Process Jira issues 2026-01-20 05:39:17 UTC
#
# For each processed Jira issue row, record the "Closed" (or "with Local
# Security") timestamp and the computed duration-in-hours metric that was
# produced by the processing run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; DateCol="D"; DateVal="2026-01-19T21:22:53.253-0500"; NumCol="H"; NumVal="5.43"}
    @{Row=3; DateCol="D"; DateVal="2026-01-19T21:23:05.223-0500"; NumCol="H"; NumVal="5.42"}
    @{Row=4; DateCol="D"; DateVal="2026-01-19T21:23:16.810-0500"; NumCol="H"; NumVal="5.42"}
    @{Row=5; DateCol="D"; DateVal="2026-01-19T21:23:28.997-0500"; NumCol="H"; NumVal="5.43"}
    @{Row=6; DateCol="D"; DateVal="2026-01-19T21:23:40.830-0500"; NumCol="H"; NumVal="5.43"}
    @{Row=7; DateCol="D"; DateVal="2026-01-19T21:23:46.903-0500"; NumCol="H"; NumVal="5.43"}
    @{Row=8; DateCol="D"; DateVal="2026-01-19T21:24:02.198-0500"; NumCol="H"; NumVal="5.43"}
    @{Row=9; DateCol="D"; DateVal="2026-01-19T21:24:09.609-0500"; NumCol="H"; NumVal="5.43"}
    @{Row=10; DateCol="D"; DateVal="2026-01-19T21:24:16.787-0500"; NumCol="H"; NumVal="5.44"}
    @{Row=11; DateCol="D"; DateVal="2026-01-19T21:24:31.814-0500"; NumCol="H"; NumVal="5.44"}
    @{Row=12; DateCol="D"; DateVal="2026-01-19T21:24:46.134-0500"; NumCol="H"; NumVal="5.44"}
    @{Row=13; DateCol="D"; DateVal="2026-01-19T21:24:58.266-0500"; NumCol="H"; NumVal="5.45"}
    @{Row=14; DateCol="D"; DateVal="2026-01-19T21:25:13.770-0500"; NumCol="H"; NumVal="5.45"}
    @{Row=15; DateCol="D"; DateVal="2026-01-19T21:25:51.286-0500"; NumCol="H"; NumVal="5.46"}
    @{Row=16; DateCol="D"; DateVal="2026-01-19T21:26:15.859-0500"; NumCol="H"; NumVal="5.46"}
    @{Row=17; DateCol="D"; DateVal="2026-01-19T21:26:32.391-0500"; NumCol="H"; NumVal="5.47"}
    @{Row=18; DateCol="D"; DateVal="2026-01-19T21:26:46.389-0500"; NumCol="H"; NumVal="5.47"}
    @{Row=19; DateCol="D"; DateVal="2026-01-19T21:27:12.819-0500"; NumCol="H"; NumVal="5.48"}
    @{Row=20; DateCol="D"; DateVal="2026-01-19T21:27:30.403-0500"; NumCol="H"; NumVal="5.48"}
    @{Row=21; DateCol="D"; DateVal="2026-01-19T21:27:38.931-0500"; NumCol="H"; NumVal="5.48"}
    @{Row=22; DateCol="D"; DateVal="2026-01-19T21:28:51.616-0500"; NumCol="H"; NumVal="5.50"}
    @{Row=23; DateCol="D"; DateVal="2026-01-19T21:29:23.985-0500"; NumCol="H"; NumVal="5.50"}
    @{Row=24; DateCol="D"; DateVal="2026-01-19T21:29:36.050-0500"; NumCol="H"; NumVal="5.51"}
    @{Row=25; DateCol="D"; DateVal="2026-01-19T21:29:38.558-0500"; NumCol="H"; NumVal="5.51"}
    @{Row=26; DateCol="D"; DateVal="2026-01-19T21:29:45.748-0500"; NumCol="H"; NumVal="5.51"}
    @{Row=27; DateCol="D"; DateVal="2026-01-19T21:29:55.804-0500"; NumCol="H"; NumVal="5.51"}
    @{Row=28; DateCol="D"; DateVal="2026-01-19T21:30:09.034-0500"; NumCol="H"; NumVal="5.51"}
    @{Row=29; DateCol="D"; DateVal="2026-01-19T21:30:23.329-0500"; NumCol="H"; NumVal="5.51"}
    @{Row=30; DateCol="D"; DateVal="2026-01-19T21:30:57.189-0500"; NumCol="H"; NumVal="5.52"}
    @{Row=31; DateCol="D"; DateVal="2026-01-19T21:32:12.560-0500"; NumCol="H"; NumVal="5.54"}
    @{Row=32; DateCol="D"; DateVal="2026-01-19T21:32:29.835-0500"; NumCol="H"; NumVal="5.55"}
    @{Row=33; DateCol="D"; DateVal="2026-01-19T21:32:54.406-0500"; NumCol="H"; NumVal="5.57"}
    @{Row=34; DateCol="D"; DateVal="2026-01-19T21:33:08.332-0500"; NumCol="H"; NumVal="5.57"}
    @{Row=35; DateCol="D"; DateVal="2026-01-19T21:35:40.389-0500"; NumCol="H"; NumVal="5.61"}
    @{Row=36; DateCol="D"; DateVal="2026-01-19T22:00:00.090-0500"; NumCol="H"; NumVal="6.01"}
    @{Row=37; DateCol="D"; DateVal="2026-01-19T22:00:14.281-0500"; NumCol="H"; NumVal="6.02"}
    @{Row=38; DateCol="D"; DateVal="2026-01-19T21:36:36.268-0500"; NumCol="H"; NumVal="5.63"}
    @{Row=39; DateCol="D"; DateVal="2026-01-19T21:36:58.388-0500"; NumCol="H"; NumVal="5.64"}
    @{Row=40; DateCol="D"; DateVal="2026-01-19T21:37:03.762-0500"; NumCol="H"; NumVal="5.67"}
    @{Row=41; DateCol="D"; DateVal="2026-01-19T21:37:08.979-0500"; NumCol="H"; NumVal="5.67"}
    @{Row=42; DateCol="D"; DateVal="2026-01-19T21:38:13.654-0500"; NumCol="H"; NumVal="6.06"}
    @{Row=45; DateCol="C"; DateVal="2026-01-19T22:38:37.999-0500"; NumCol="G"; NumVal="10.88"}
    @{Row=51; DateCol="C"; DateVal="2026-01-19T21:27:35.536-0500"; NumCol="G"; NumVal="7.76"}
)

foreach ($u in $updates) {
    $dateCell = $ws.Range("$($u.DateCol)$($u.Row)")
    $dateCell.Value = $u.DateVal

    # The metric column holds a numeric-looking string ("5.43", "10.88", ...)
    # that must stay plain text (matches the sheet's existing inline-string
    # cells), so force text formatting before assigning it.
    $numCell = $ws.Range("$($u.NumCol)$($u.Row)")
    $numCell.NumberFormat = "@"
    $numCell.Value = $u.NumVal
}
